$wb = $excel.ActiveWorkbook
$lastIdx = $wb.Worksheets.Count
$last = $wb.Worksheets.Item($lastIdx)
$ws = $wb.Worksheets.Add($null, $last)
$ws.Name = "L6"

$headers = @("Form", "Goals scored", "Goals conceded", "Total Goals")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$ws.Range("A2:A25").NumberFormat = "@"

$rows = @(
    @("1", "Barnsley,W W L W W L", "Barnsley,2 2 0 1 1 0", "Barnsley,1 0 2 0 0 2", "Barnsley,3 2 2 1 1 2"),
    @("2", "Birmingham,D W W D W L", "Birmingham,0 2 1 1 2 0", "Birmingham,0 0 0 1 1 4", "Birmingham,0 2 1 2 3 4"),
    @("3", "Blackburn,L D W L W D", "Blackburn,0 2 2 0 5 1", "Blackburn,2 2 1 1 2 1", "Blackburn,2 4 3 1 7 2"),
    @("4", "Bournemouth,W W W W L L", "Bournemouth,4 2 3 4 0 0", "Bournemouth,1 1 1 1 1 1", "Bournemouth,5 3 4 5 1 1"),
    @("5", "Brentford,W D D W W W", "Brentford,5 0 1 1 1 2", "Brentford,0 0 1 0 0 0", "Brentford,5 0 2 1 1 2"),
    @("6", "Bristol City,L D D L L L", "Bristol City,1 0 1 1 2 1", "Bristol City,3 0 1 2 3 4", "Bristol City,4 0 2 3 5 5"),
    @("7", "Cardiff,L D D D W W", "Cardiff,0 2 1 1 2 4", "Cardiff,5 2 1 1 1 0", "Cardiff,5 4 2 2 3 4"),
    @("8", "Coventry,L W W W L D", "Coventry,1 1 2 3 0 1", "Coventry,4 0 0 2 1 1", "Coventry,5 1 2 5 1 2"),
    @("9", "Derby,L L L L L L", "Derby,1 0 1 0 1 1", "Derby,3 1 2 3 2 2", "Derby,4 1 3 3 3 3"),
    @("10", "Huddersfield,D L W L L D", "Huddersfield,0 1 2 0 2 1", "Huddersfield,0 2 0 1 5 1", "Huddersfield,0 3 2 1 7 2"),
    @("11", "Luton,L W W D W D", "Luton,1 3 1 0 3 1", "Luton,2 1 0 0 2 1", "Luton,3 4 1 0 5 2"),
    @("12", "Middlesbrough,D L L W W D", "Middlesbrough,1 0 1 2 3 1", "Middlesbrough,1 2 2 1 1 1", "Middlesbrough,2 2 3 3 4 2"),
    @("13", "Millwall,W L D L L W", "Millwall,2 0 0 1 0 4", "Millwall,1 3 0 4 1 1", "Millwall,3 3 0 5 1 5"),
    @("14", "Norwich,W W L L W W", "Norwich,7 1 1 0 3 4", "Norwich,0 0 3 1 1 1", "Norwich,7 1 4 1 4 5"),
    @("15", "Nottm Forest,W D L D D D", "Nottm Forest,3 0 0 1 1 0", "Nottm Forest,1 0 2 1 1 0", "Nottm Forest,4 0 2 2 2 0"),
    @("16", "Preston,W L D W W W", "Preston,1 0 0 3 1 2", "Preston,0 5 0 0 0 0", "Preston,1 5 0 3 1 2"),
    @("17", "QPR,W L W W L W", "QPR,4 1 2 1 1 2", "QPR,1 3 1 0 3 0", "QPR,5 4 3 1 4 2"),
    @("18", "Reading,W L D D D L", "Reading,3 0 1 0 2 1", "Reading,1 2 1 0 2 4", "Reading,4 2 2 0 4 5"),
    @("19", "Rotherham,L L L L L D", "Rotherham,0 0 1 0 0 1", "Rotherham,1 1 2 1 1 1", "Rotherham,1 1 3 1 1 2"),
    @("20", "Sheffield Weds,L L D W L D", "Sheffield Weds,1 0 1 1 1 0", "Sheffield Weds,4 2 1 0 3 0", "Sheffield Weds,5 2 2 1 4 0"),
    @("21", "Stoke,L L D L D L", "Stoke,1 0 0 2 1 0", "Stoke,2 2 0 3 1 2", "Stoke,3 2 0 5 2 2"),
    @("22", "Swansea,W W D L D W", "Swansea,3 2 2 0 2 2", "Swansea,0 0 2 1 2 1", "Swansea,3 2 4 1 4 3"),
    @("23", "Watford,D W L W W L", "Watford,1 2 0 1 1 0", "Watford,1 0 1 0 0 2", "Watford,2 2 1 1 1 2"),
    @("24", "Wycombe,W L D W L W", "Wycombe,3 1 2 2 1 1", "Wycombe,0 3 2 1 2 0", "Wycombe,3 4 4 3 3 1")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
